$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Benchmarks")

# --- Selection: place the active cell at I31 (row 20 area work) ---
[void]$ws.Range("I31").Select()

# --- Row 20: new Array Dimension benchmark entry ---
$ws.Range("F20").Value = 4
$ws.Range("G19").Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("G20").Value = "'(5, 5, 7, 7, 7)"
$ws.Range("I20").Formula = "=5.593/4"
$ws.Range("K20").Formula = "=I20/`$I`$13"

# --- Row 23: updated Gaussian filter timings ---
$ws.Range("I23").Formula = "=19.617/4"

# --- Row 25 ---
$ws.Range("I25").Formula = "=25.726/4"

# --- Row 26 ---
$ws.Range("I26").Formula = "=2.005/4"

# --- Row 27: mark as Production ---
$ws.Range("D27").Value = "Production "

# --- Row 28: mark as Production + updated timing ---
$ws.Range("D28").Value = "Production "
$ws.Range("I28").Formula = "=0.731/4"

# --- Row 29: mark as Production + updated timing ---
$ws.Range("D29").Value = "Production "
$ws.Range("I29").Formula = "=0.164/4"

# --- Row 30: mark as Production + new benchmark entry ---
$ws.Range("D30").Value = "Production "
$ws.Range("F30").Value = 4
$ws.Range("G19").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = "'(5, 5, 7, 7, 7)"
$ws.Range("I30").Formula = "=0.114/4"
$ws.Range("K30").Formula = "=I30/`$I`$22"

# --- Row 38: kernel label corrected to (5, 5, 7, 7, 7) ---
$ws.Range("G38").Value = "'(5, 5, 7, 7, 7)"

# --- Row 39: new benchmark entry ---
$ws.Range("F39").Value = 4
$ws.Range("G19").Copy()
$ws.Range("G39").PasteSpecial(-4122)
$ws.Range("G39").Value = "'(5, 5, 7, 7, 7)"
$ws.Range("K38").Copy()
$ws.Range("K39").PasteSpecial(-4122)
$ws.Range("I39").Value = 0.067
$ws.Range("K39").Formula = "=I39/`$I`$32"

[void]($excel.CutCopyMode = 0)
